$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text (shared string used by A1) to the corrected Kyrgyz wording
$ws.Range("A1").Value = "8.10.2.2 Камсыздандыруу компаниялардын финансылык көрсөткүчтөрү"

# Add a new column Q (year 2023) mirroring the existing column P values/styles
$ws.Range("Q3").Value = 2023
$ws.Range("Q4").Value = 16
$ws.Range("Q5").Value = 3031.4

$ws.Range("P3:P5").Copy()
$ws.Range("Q3:Q5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reset the cursor to A1 (the stray "E9" selection left in the sheet view is stale)
$ws.Range("A1").Select()
